# Applies the weekly Fruta/Hortaliza update to the Repollo (Terminal Hortofrutícola
# Agro Chillán) subset sheet: rows 537-546 get revised figures, and the historical
# record that used to sit in row 547 is expanded back out into four distinct rows
# (547-550), each keeping the shared market/category columns but with its own date,
# volume and price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 537-546 (field-level changes) ---
$ws.Cells.Item(537, 4).Value = 45239
$ws.Cells.Item(537, 10).Value = 500

$ws.Cells.Item(538, 4).Value = 45239
$ws.Cells.Item(538, 10).Value = 400
$ws.Cells.Item(538, 11).Value = 1200
$ws.Cells.Item(538, 12).Value = 1200
$ws.Cells.Item(538, 13).Value = 1200
$ws.Cells.Item(538, 15).Value = "Región del Maule"
$ws.Cells.Item(538, 16).Value = 1200

$ws.Cells.Item(539, 4).Value = 45239
$ws.Cells.Item(539, 9).Value = "Segunda"
$ws.Cells.Item(539, 10).Value = 500
$ws.Cells.Item(539, 11).Value = 1000
$ws.Cells.Item(539, 12).Value = 1000
$ws.Cells.Item(539, 13).Value = 1000
$ws.Cells.Item(539, 15).Value = "Región del Maule"
$ws.Cells.Item(539, 16).Value = 1000

$ws.Cells.Item(540, 4).Value = 45173
$ws.Cells.Item(540, 9).Value = "Primera"
$ws.Cells.Item(540, 10).Value = 250
$ws.Cells.Item(540, 11).Value = 1000
$ws.Cells.Item(540, 12).Value = 1000
$ws.Cells.Item(540, 13).Value = 1000
$ws.Cells.Item(540, 16).Value = 1000

$ws.Cells.Item(541, 4).Value = 44607
$ws.Cells.Item(541, 10).Value = 200
$ws.Cells.Item(541, 11).Value = 700
$ws.Cells.Item(541, 12).Value = 750
$ws.Cells.Item(541, 13).Value = 725
$ws.Cells.Item(541, 16).Value = 725

$ws.Cells.Item(542, 4).Value = 44858
$ws.Cells.Item(542, 9).Value = "Primera"
$ws.Cells.Item(542, 10).Value = 400
$ws.Cells.Item(542, 11).Value = 1500
$ws.Cells.Item(542, 12).Value = 1600
$ws.Cells.Item(542, 13).Value = 1550
$ws.Cells.Item(542, 16).Value = 1550

$ws.Cells.Item(543, 4).Value = 44858
$ws.Cells.Item(543, 9).Value = "Segunda"
$ws.Cells.Item(543, 10).Value = 300
$ws.Cells.Item(543, 11).Value = 1200
$ws.Cells.Item(543, 12).Value = 1200
$ws.Cells.Item(543, 13).Value = 1200
$ws.Cells.Item(543, 16).Value = 1200

$ws.Cells.Item(544, 4).Value = 45077
$ws.Cells.Item(544, 10).Value = 500
$ws.Cells.Item(544, 11).Value = 1200
$ws.Cells.Item(544, 12).Value = 1300
$ws.Cells.Item(544, 13).Value = 1250
$ws.Cells.Item(544, 16).Value = 1250

$ws.Cells.Item(545, 4).Value = 45077
$ws.Cells.Item(545, 9).Value = "Segunda"
$ws.Cells.Item(545, 10).Value = 300
$ws.Cells.Item(545, 11).Value = 1000
$ws.Cells.Item(545, 12).Value = 1000
$ws.Cells.Item(545, 13).Value = 1000
$ws.Cells.Item(545, 16).Value = 1000

$ws.Cells.Item(546, 4).Value = 44386
$ws.Cells.Item(546, 10).Value = 160
$ws.Cells.Item(546, 11).Value = 700
$ws.Cells.Item(546, 12).Value = 750
$ws.Cells.Item(546, 13).Value = 725
$ws.Cells.Item(546, 16).Value = 725

# --- Row 547 historical entry expands into rows 547-550 ---
$ws.Cells.Item(547, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(547, 1).Value = 7
$ws.Cells.Item(547, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(547, 3).Value = "Ñuble"
$ws.Cells.Item(547, 4).Value = 44477
$ws.Cells.Item(547, 5).Value = 16
$ws.Cells.Item(547, 6).Value = 100112006
$ws.Cells.Item(547, 7).Value = "Repollo"
$ws.Cells.Item(547, 8).Value = "Crespo record"
$ws.Cells.Item(547, 9).Value = "Primera"
$ws.Cells.Item(547, 10).Value = 300
$ws.Cells.Item(547, 11).Value = 600
$ws.Cells.Item(547, 12).Value = 650
$ws.Cells.Item(547, 13).Value = 625
$ws.Cells.Item(547, 14).Value = "$/unidad"
$ws.Cells.Item(547, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(547, 16).Value = 625
$ws.Cells.Item(547, 17).Value = 1
$ws.Cells.Item(547, 18).Value = "Hortaliza"

$ws.Cells.Item(548, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(548, 1).Value = 7
$ws.Cells.Item(548, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(548, 3).Value = "Ñuble"
$ws.Cells.Item(548, 4).Value = 44508
$ws.Cells.Item(548, 5).Value = 16
$ws.Cells.Item(548, 6).Value = 100112006
$ws.Cells.Item(548, 7).Value = "Repollo"
$ws.Cells.Item(548, 8).Value = "Crespo record"
$ws.Cells.Item(548, 9).Value = "Primera"
$ws.Cells.Item(548, 10).Value = 400
$ws.Cells.Item(548, 11).Value = 600
$ws.Cells.Item(548, 12).Value = 700
$ws.Cells.Item(548, 13).Value = 650
$ws.Cells.Item(548, 14).Value = "$/unidad"
$ws.Cells.Item(548, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(548, 16).Value = 650
$ws.Cells.Item(548, 17).Value = 1
$ws.Cells.Item(548, 18).Value = "Hortaliza"

$ws.Cells.Item(549, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(549, 1).Value = 7
$ws.Cells.Item(549, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(549, 3).Value = "Ñuble"
$ws.Cells.Item(549, 4).Value = 44579
$ws.Cells.Item(549, 5).Value = 16
$ws.Cells.Item(549, 6).Value = 100112006
$ws.Cells.Item(549, 7).Value = "Repollo"
$ws.Cells.Item(549, 8).Value = "Crespo record"
$ws.Cells.Item(549, 9).Value = "Primera"
$ws.Cells.Item(549, 10).Value = 300
$ws.Cells.Item(549, 11).Value = 600
$ws.Cells.Item(549, 12).Value = 700
$ws.Cells.Item(549, 13).Value = 650
$ws.Cells.Item(549, 14).Value = "$/unidad"
$ws.Cells.Item(549, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(549, 16).Value = 650
$ws.Cells.Item(549, 17).Value = 1
$ws.Cells.Item(549, 18).Value = "Hortaliza"

$ws.Cells.Item(550, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(550, 1).Value = 7
$ws.Cells.Item(550, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(550, 3).Value = "Ñuble"
$ws.Cells.Item(550, 4).Value = 44778
$ws.Cells.Item(550, 5).Value = 16
$ws.Cells.Item(550, 6).Value = 100112006
$ws.Cells.Item(550, 7).Value = "Repollo"
$ws.Cells.Item(550, 8).Value = "Crespo record"
$ws.Cells.Item(550, 9).Value = "Primera"
$ws.Cells.Item(550, 10).Value = 240
$ws.Cells.Item(550, 11).Value = 1100
$ws.Cells.Item(550, 12).Value = 1200
$ws.Cells.Item(550, 13).Value = 1150
$ws.Cells.Item(550, 14).Value = "$/unidad"
$ws.Cells.Item(550, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(550, 16).Value = 1150
$ws.Cells.Item(550, 17).Value = 1
$ws.Cells.Item(550, 18).Value = "Hortaliza"

